$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9578740567487785
$ws.Range("D2").Value = 0.2826761273864307
$ws.Range("E2").Value = 0.2790990945624525
$ws.Range("F2").Value = 1.378316428980874
$ws.Range("G2").Value = 0.7021962693790087
$ws.Range("H2").Value = 0.8157896895588266
$ws.Range("J2").Value = 0.3486991817280822
$ws.Range("K2").Value = 0.4024466136276601
$ws.Range("L2").Value = 0.1258011187399575
$ws.Range("O2").Value = 3.036201164720595
$ws.Range("B3").Value = 0.9296918391557369
$ws.Range("D3").Value = 0.2826239691642058
$ws.Range("E3").Value = 0.2812468761275255
$ws.Range("F3").Value = 1.387076401537072
$ws.Range("G3").Value = 0.7072734149181557
$ws.Range("H3").Value = 0.8219323259065874
$ws.Range("J3").Value = 0.3519234999010887
$ws.Range("K3").Value = 0.3519541583104626
$ws.Range("L3").Value = 0.1148024372222238
$ws.Range("O3").Value = 3.059571074234356
$ws.Range("B4").Value = 0.9127114802716108
$ws.Range("D4").Value = 0.2826914435864722
$ws.Range("E4").Value = 0.2826616876326433
$ws.Range("F4").Value = 1.393129313884096
$ws.Range("G4").Value = 0.7108170665987785
$ws.Range("H4").Value = 0.8260294332101523
$ws.Range("J4").Value = 0.3540163459668295
$ws.Range("K4").Value = 0.3208151669452377
$ws.Range("L4").Value = 0.1080527022151614
$ws.Range("O4").Value = 3.075496523896987
$ws.Range("B5").Value = 0.905873922664739
$ws.Range("D5").Value = 0.2827440390969898
$ws.Range("E5").Value = 0.2832624290474701
$ws.Range("F5").Value = 1.395765613410376
$ws.Range("G5").Value = 0.7123683196237991
$ws.Range("H5").Value = 0.8277809698565903
$ws.Range("J5").Value = 0.3548976750406649
$ws.Range("K5").Value = 0.3080923436759804
$ws.Range("L5").Value = 0.1053032180933684
$ws.Range("O5").Value = 3.082382791975945
$ws.Range("B6").Value = 0.9047435291374484
$ws.Range("D6").Value = 0.2827542907925746
$ws.Range("E6").Value = 0.283363644088551
$ws.Range("F6").Value = 1.396213621879248
$ws.Range("G6").Value = 0.7126323777462673
$ws.Range("H6").Value = 0.8280767620816007
$ws.Range("J6").Value = 0.3550457401090545
$ws.Range("K6").Value = 0.3059777351338084
$ws.Range("L6").Value = 0.1048467401651152
$ws.Range("O6").Value = 3.083550205011363
$ws.Range("B7").Value = 0.9126189333308616
$ws.Range("D7").Value = 0.2826920511767312
$ws.Range("E7").Value = 0.2826696914292697
$ws.Range("F7").Value = 1.393164180690384
$ws.Range("G7").Value = 0.7108375533490303
$ws.Range("H7").Value = 0.8260527231976198
$ws.Range("J7").Value = 0.3540281165256127
$ws.Range("K7").Value = 0.3206437167049501
$ws.Range("L7").Value = 0.1080156170307873
$ws.Range("O7").Value = 3.075587788816051
$ws.Range("B8").Value = 0.94809006371122
$ws.Range("D8").Value = 0.2826375389411879
$ws.Range("E8").Value = 0.2798197342383038
$ws.Range("F8").Value = 1.381197010602037
$ws.Range("G8").Value = 0.7038583968575622
$ws.Range("H8").Value = 0.8178401612967861
$ws.Range("J8").Value = 0.3497874748709644
$ws.Range("K8").Value = 0.3850656362081679
$ws.Range("L8").Value = 0.1220081678903284
$ws.Range("O8").Value = 3.043932059174892
$ws.Range("B9").Value = 1.020190384363588
$ws.Range("D9").Value = 0.2833168811184521
$ws.Range("E9").Value = 0.2749915033959116
$ws.Range("F9").Value = 1.363073586058221
$ws.Range("G9").Value = 0.6935549431757835
$ws.Range("H9").Value = 0.8043145146338375
$ws.Range("J9").Value = 0.3423672525352424
$ws.Range("K9").Value = 0.5102819989655529
$ws.Range("L9").Value = 0.1494675222748043
$ws.Range("O9").Value = 2.994355697919218
$ws.Range("B10").Value = 1.074682598149678
$ws.Range("D10").Value = 0.2842912822583372
$ws.Range("E10").Value = 0.2719054882064054
$ws.Range("F10").Value = 1.353009267767135
$ws.Range("G10").Value = 0.6880483094427348
$ws.Range("H10").Value = 0.7959447883269775
$ws.Range("J10").Value = 0.3374591027289719
$ws.Range("K10").Value = 0.6015640175356509
$ws.Range("L10").Value = 0.1696456637964587
$ws.Range("O10").Value = 2.965545417259193
$ws.Range("B11").Value = 1.099796775334767
$ws.Range("D11").Value = 0.284836935830981
$ws.Range("E11").Value = 0.2706012450835082
$ws.Range("F11").Value = 1.349135238458459
$ws.Range("G11").Value = 0.6859915113777362
$ws.Range("H11").Value = 0.7924765698010248
$ws.Range("J11").Value = 0.3353437298418438
$ws.Range("K11").Value = 0.6429286281858424
$ws.Range("L11").Value = 0.1788243754380119
$ws.Range("O11").Value = 2.954090665783639
$ws.Range("B12").Value = 1.109353032087625
$ws.Range("D12").Value = 0.2850582094049514
$ws.Range("E12").Value = 0.2701216463928056
$ws.Range("F12").Value = 1.347769392803059
$ws.Range("G12").Value = 0.6852771209645283
$ws.Range("H12").Value = 0.7912119466888612
$ws.Range("J12").Value = 0.33455953266101
$ws.Range("K12").Value = 0.65856855410496
$ws.Range("L12").Value = 0.1822998699663998
$ws.Range("O12").Value = 2.949990365330848
$ws.Range("B13").Value = 1.107292884140008
$ws.Range("D13").Value = 0.2850099040156806
$ws.Range("E13").Value = 0.2702243015239585
$ws.Range("F13").Value = 1.348059054505278
$ws.Range("G13").Value = 0.685428109832813
$ws.Range("H13").Value = 0.791482140237008
$ws.Range("J13").Value = 0.3347276746651886
$ws.Range("K13").Value = 0.6552012965445897
$ws.Range("L13").Value = 0.1815513762273611
$ws.Range("O13").Value = 2.95086288270403
$ws.Range("B14").Value = 1.100582055329028
$ws.Range("D14").Value = 0.2848548470499139
$ws.Range("E14").Value = 0.2705615020066237
$ws.Range("F14").Value = 1.349020842642908
$ws.Range("G14").Value = 0.6859314457510806
$ws.Range("H14").Value = 0.7923715525115114
$ws.Range("J14").Value = 0.3352788759608023
$ws.Range("K14").Value = 0.6442158197625361
$ws.Range("L14").Value = 0.1791103134390255
$ws.Range("O14").Value = 2.953748574871554
$ws.Range("B15").Value = 1.096477456058494
$ws.Range("D15").Value = 0.2847617752725142
$ws.Range("E15").Value = 0.2707699070367564
$ws.Range("F15").Value = 1.349623137080364
$ws.Range("G15").Value = 0.6862481507015019
$ws.Range("H15").Value = 0.792922685623509
$ws.Range("J15").Value = 0.3356186960043122
$ws.Range("K15").Value = 0.6374837506013762
$ws.Range("L15").Value = 0.1776150483544399
$ws.Range("O15").Value = 2.955547054401279
$ws.Range("B16").Value = 1.073047813087442
$ws.Range("D16").Value = 0.2842576762949562
$ws.Range("E16").Value = 0.2719927249946679
$ws.Range("F16").Value = 1.353276607265876
$ws.Range("G16").Value = 0.6881917469954004
$ws.Range("H16").Value = 0.7961782650466205
$ws.Range("J16").Value = 0.3375997066106811
$ws.Range("K16").Value = 0.5988574420321982
$ws.Range("L16").Value = 0.1690457857431227
$ws.Range("O16").Value = 2.966327232102003
$ws.Range("B17").Value = 1.058757308828149
$ws.Range("D17").Value = 0.2839745980763126
$ws.Range("E17").Value = 0.2727683711114501
$ws.Range("F17").Value = 1.355698198190993
$ws.Range("G17").Value = 0.6894988919331411
$ws.Range("H17").Value = 0.7982622935145258
$ws.Range("J17").Value = 0.3388450328356609
$ws.Range("K17").Value = 0.5751197923684686
$ws.Range("L17").Value = 0.1637885525943119
$ws.Range("O17").Value = 2.973363379331772
$ws.Range("B18").Value = 1.050568460011306
$ws.Range("D18").Value = 0.2838214213843031
$ws.Range("E18").Value = 0.2732238790714145
$ws.Range("F18").Value = 1.35715732995375
$ws.Range("G18").Value = 0.6902929121911612
$ws.Range("H18").Value = 0.7994929035135172
$ws.Range("J18").Value = 0.3395723605882672
$ws.Range("K18").Value = 0.5614515166295462
$ws.Range("L18").Value = 0.160764704317927
$ws.Range("O18").Value = 2.977565807423034
$ws.Range("B19").Value = 1.04780114654892
$ws.Range("D19").Value = 0.2837712168485353
$ws.Range("E19").Value = 0.2733797177134711
$ws.Range("F19").Value = 1.357662756742997
$ws.Range("G19").Value = 0.6905689983811882
$ws.Range("H19").Value = 0.7999150531589336
$ws.Range("J19").Value = 0.3398205198371671
$ws.Range("K19").Value = 0.5568211245153805
$ws.Range("L19").Value = 0.1597408834508656
$ws.Range("O19").Value = 2.979015373519303
$ws.Range("B20").Value = 1.060275388324044
$ws.Range("D20").Value = 0.2840037348456974
$ws.Range("E20").Value = 0.2726848320282329
$ws.Range("F20").Value = 1.355433555032128
$ws.Range("G20").Value = 0.6893553779490844
$ws.Range("H20").Value = 0.7980371407406892
$ws.Range("J20").Value = 0.3387113224090239
$ws.Range("K20").Value = 0.5776482657489339
$ws.Range("L20").Value = 0.1643481984324495
$ws.Range("O20").Value = 2.972598284349772
$ws.Range("B21").Value = 1.102551944089782
$ws.Range("D21").Value = 0.2848999941589341
$ws.Range("E21").Value = 0.2704620704870369
$ws.Range("F21").Value = 1.348735597434725
$ws.Range("G21").Value = 0.6857818537872618
$ws.Range("H21").Value = 0.7921089889250652
$ws.Range("J21").Value = 0.3351165177413273
$ws.Range("K21").Value = 0.6474431781809074
$ws.Range("L21").Value = 0.1798273222442788
$ws.Range("O21").Value = 2.952894535701205
$ws.Range("B22").Value = 1.130450269848296
$ws.Range("D22").Value = 0.2855710859031859
$ws.Range("E22").Value = 0.2690926445982615
$ws.Range("F22").Value = 1.344947711672013
$ws.Range("G22").Value = 0.6838221656891932
$ws.Range("H22").Value = 0.788518524213714
$ws.Range("J22").Value = 0.3328652985149976
$ws.Range("K22").Value = 0.6929181297251716
$ws.Range("L22").Value = 0.1899420639228424
$ws.Range("O22").Value = 2.941400508261694
$ws.Range("B23").Value = 1.115536112760765
$ws.Range("D23").Value = 0.2852051291450692
$ws.Range("E23").Value = 0.269815923830814
$ws.Range("F23").Value = 1.346915463954247
$ws.Range("G23").Value = 0.6848336931438013
$ws.Range("H23").Value = 0.7904088641643199
$ws.Range("J23").Value = 0.3340578413616502
$ws.Range("K23").Value = 0.6686604175088462
$ws.Range("L23").Value = 0.1845438705540801
$ws.Range("O23").Value = 2.947408520462488
$ws.Range("B24").Value = 1.059588980575853
$ws.Range("D24").Value = 0.283990532301587
$ws.Range("E24").Value = 0.2727225702090603
$ws.Range("F24").Value = 1.355552991724593
$ws.Range("G24").Value = 0.6894201281605845
$ws.Range("H24").Value = 0.7981388311499131
$ws.Range("J24").Value = 0.3387717374648211
$ws.Range("K24").Value = 0.5765052088511879
$ws.Range("L24").Value = 0.1640951868781002
$ws.Range("O24").Value = 2.972943693952914
$ws.Range("B25").Value = 1.000416280950645
$ws.Range("D25").Value = 0.2830493449624711
$ws.Range("E25").Value = 0.2762164916344023
$ws.Range("F25").Value = 1.367405043129303
$ws.Range("G25").Value = 0.6959800414851998
$ws.Range("H25").Value = 0.8076979363624162
$ws.Range("J25").Value = 0.3442789952225755
$ws.Range("K25").Value = 0.4765306208675213
$ws.Range("L25").Value = 0.1420377731562326
$ws.Range("O25").Value = 3.00642989212227
